# Fruta / hortaliza, semanal
# Insert two new price-observation rows (new rows 639 and 640) into the
# "Vega Modelo de Temuco - Uva" sheet, pushing the existing rows 639-684
# down to 641-686.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 639, shifting everything from old row 639 onward
# down by two rows (old 639 -> new 641, ... old 684 -> new 686).
$ws.Rows("639:640").Insert()

# --- New row 639 ---
$ws.Cells.Item(639, 1).Value  = 10
$ws.Cells.Item(639, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(639, 3).Value  = "La Araucanía"
$ws.Cells.Item(639, 4).Value  = 44610
$ws.Cells.Item(639, 5).Value  = 9
$ws.Cells.Item(639, 6).Value  = "Fruta"
$ws.Cells.Item(639, 7).Value  = 100109
$ws.Cells.Item(639, 8).Value  = "Uva"
$ws.Cells.Item(639, 9).Value  = 100109001
$ws.Cells.Item(639, 10).Value = "Uva"
$ws.Cells.Item(639, 11).Value = "Red Globe"
$ws.Cells.Item(639, 12).Value = "Primera"
$ws.Cells.Item(639, 13).Value = 300
$ws.Cells.Item(639, 14).Value = 13000
$ws.Cells.Item(639, 15).Value = 13000
$ws.Cells.Item(639, 16).Value = 13000
$ws.Cells.Item(639, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(639, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(639, 19).Value = 722
$ws.Cells.Item(639, 20).Value = 18

# --- New row 640 ---
$ws.Cells.Item(640, 1).Value  = 10
$ws.Cells.Item(640, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(640, 3).Value  = "La Araucanía"
$ws.Cells.Item(640, 4).Value  = 44610
$ws.Cells.Item(640, 5).Value  = 9
$ws.Cells.Item(640, 6).Value  = "Fruta"
$ws.Cells.Item(640, 7).Value  = 100109
$ws.Cells.Item(640, 8).Value  = "Uva"
$ws.Cells.Item(640, 9).Value  = 100109001
$ws.Cells.Item(640, 10).Value = "Uva"
$ws.Cells.Item(640, 11).Value = "Thompson seedless"
$ws.Cells.Item(640, 12).Value = "Primera"
$ws.Cells.Item(640, 13).Value = 200
$ws.Cells.Item(640, 14).Value = 13000
$ws.Cells.Item(640, 15).Value = 13000
$ws.Cells.Item(640, 16).Value = 13000
$ws.Cells.Item(640, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(640, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(640, 19).Value = 722
$ws.Cells.Item(640, 20).Value = 18

# Ensure the date cells keep the workbook's date number format (style index 2,
# same as every other "Fecha" cell in column D).
$ws.Cells.Item(639, 4).NumberFormat = $ws.Cells.Item(641, 4).NumberFormat
$ws.Cells.Item(640, 4).NumberFormat = $ws.Cells.Item(641, 4).NumberFormat
